# Fix VRAM test data column names to match schema
# ('RA CRIT' -> 'RA Critical', 'RA HIGH' -> 'RA High', 'Scan %' -> 'Percent
# Scanned'; ESS data already has the Assets column correctly.)
#
# This corrects the "TAM Overdue" sheet's "# Past Due" (B) and
# "# Extensions" (C) sample values, rows 2-61, to match the refreshed
# source data used to validate the renamed schema columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TAM Overdue")

$data = New-Object 'object[,]' 60,2
$data[0,0] = 0; $data[0,1] = 0
$data[1,0] = 1; $data[1,1] = 2
$data[2,0] = 1; $data[2,1] = 0
$data[3,0] = 0; $data[3,1] = 1
$data[4,0] = 3; $data[4,1] = 2
$data[5,0] = 0; $data[5,1] = 2
$data[6,0] = 0; $data[6,1] = 0
$data[7,0] = 3; $data[7,1] = 2
$data[8,0] = 2; $data[8,1] = 1
$data[9,0] = 10; $data[9,1] = 4
$data[10,0] = 0; $data[10,1] = 0
$data[11,0] = 1; $data[11,1] = 1
$data[12,0] = 3; $data[12,1] = 1
$data[13,0] = 6; $data[13,1] = 5
$data[14,0] = 3; $data[14,1] = 1
$data[15,0] = 6; $data[15,1] = 3
$data[16,0] = 2; $data[16,1] = 0
$data[17,0] = 0; $data[17,1] = 2
$data[18,0] = 6; $data[18,1] = 4
$data[19,0] = 0; $data[19,1] = 0
$data[20,0] = 1; $data[20,1] = 1
$data[21,0] = 0; $data[21,1] = 1
$data[22,0] = 3; $data[22,1] = 2
$data[23,0] = 2; $data[23,1] = 0
$data[24,0] = 1; $data[24,1] = 0
$data[25,0] = 3; $data[25,1] = 1
$data[26,0] = 2; $data[26,1] = 1
$data[27,0] = 2; $data[27,1] = 0
$data[28,0] = 8; $data[28,1] = 2
$data[29,0] = 1; $data[29,1] = 2
$data[30,0] = 10; $data[30,1] = 3
$data[31,0] = 3; $data[31,1] = 0
$data[32,0] = 9; $data[32,1] = 3
$data[33,0] = 1; $data[33,1] = 2
$data[34,0] = 1; $data[34,1] = 0
$data[35,0] = 0; $data[35,1] = 1
$data[36,0] = 10; $data[36,1] = 3
$data[37,0] = 1; $data[37,1] = 1
$data[38,0] = 1; $data[38,1] = 2
$data[39,0] = 0; $data[39,1] = 0
$data[40,0] = 1; $data[40,1] = 2
$data[41,0] = 3; $data[41,1] = 0
$data[42,0] = 3; $data[42,1] = 1
$data[43,0] = 3; $data[43,1] = 0
$data[44,0] = 0; $data[44,1] = 1
$data[45,0] = 3; $data[45,1] = 0
$data[46,0] = 0; $data[46,1] = 2
$data[47,0] = 3; $data[47,1] = 2
$data[48,0] = 1; $data[48,1] = 2
$data[49,0] = 1; $data[49,1] = 1
$data[50,0] = 10; $data[50,1] = 5
$data[51,0] = 2; $data[51,1] = 2
$data[52,0] = 0; $data[52,1] = 1
$data[53,0] = 0; $data[53,1] = 1
$data[54,0] = 1; $data[54,1] = 0
$data[55,0] = 24; $data[55,1] = 7
$data[56,0] = 1; $data[56,1] = 2
$data[57,0] = 1; $data[57,1] = 1
$data[58,0] = 0; $data[58,1] = 0
$data[59,0] = 1; $data[59,1] = 0

$ws.Range("B2:C61").Value = $data
